# Add the 4th new parish ("Moudon - Syens") as row 7 of the id sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 8040000000
$ws.Range("C7").Value = 8040
$ws.Range("D7").Value = 8000
$ws.Range("E7").Value = "Moudon – Syens"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "P"

# Match the author's final selection (moved from E6 down to the new E7 row).
$ws.Range("E7").Select()
